$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at 171; this shifts every existing row (171..206)
# down by one (171->172, ..., 206->207) and carries formatting (e.g. the
# date-numFmt style on column D) down with it, exactly matching the target
# sheet's <dimension ref="A1:R207"/>.
$ws.Rows.Item(171).Insert()

# Populate the newly inserted row 171 with the new weekly price-report entry.
$ws.Range("A171").Value = 7
$ws.Range("B171").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C171").Value = "Ñuble"
$ws.Range("D171").Value = (Get-Date -Year 2022 -Month 5 -Day 13 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E171").Value = 16
$ws.Range("F171").Value = 100112017
$ws.Range("G171").Value = "Apio"
$ws.Range("H171").Value = "Americana (o)"
$ws.Range("I171").Value = "Primera"
$ws.Range("J171").Value = 120
$ws.Range("K171").Value = 8500
$ws.Range("L171").Value = 9000
$ws.Range("M171").Value = 8750
$ws.Range("N171").Value = "`$/docena de matas"
$ws.Range("O171").Value = "Provincia del Elquí"
$ws.Range("P171").Value = 1458
$ws.Range("Q171").Value = 6
$ws.Range("R171").Value = "Hortaliza"
